$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: the createCustomer_ID test row loses its generated-suffix value,
#     reverting D3 back to the plain base value "Rock". ---
$ws.Cells.Item(3, 4).Value2 = "Rock"

# --- Row 4: new createSalesTeam_ID test row (was entirely empty before). ---
$ws.Cells.Item(4, 1).Value2 = "createSalesTeam_ID"
$ws.Cells.Item(4, 2).Value2 = "cool.sanu09@gmail.com"
$ws.Cells.Item(4, 3).Value2 = "sanu123456"
$ws.Cells.Item(4, 4).Value2 = "Meeting-457-544-510-346-484-30"
$ws.Cells.Item(4, 5).Value2 = "jacson-201-605-714-949-76-152"
$ws.Cells.Item(4, 6).Value2 = "jacson@gmail.com-930-27-408-594-681-945"
$ws.Cells.Item(4, 7).Value2 = "jason"

# Copy the formatting from row 3 onto row 4 (A:G) so the new row matches the
# established "data row" look (bordered cells, hyperlink style on column B).
$ws.Range("A3:G3").Copy()
$ws.Range("A4:G4").PasteSpecial(-4122)

# Add the new hyperlink for the sales-team email (auto-picks up the text
# already sitting in F4 as the display text).
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:jacson@gmail.com") | Out-Null

# Column F needs the hyperlink look too (border + Hyperlink font), same as
# column B on data rows - copy that formatting onto F4 specifically (the
# hyperlink insert above can restyle the cell, so do this last).
$ws.Range("B3").Copy()
$ws.Range("F4").PasteSpecial(-4122)

# New column F is wide enough to show the generated e-mail value.
$ws.Columns.Item(6).ColumnWidth = 16.43

# Selection moves on to the next empty row, ready for the next test case.
$ws.Range("D5").Select() | Out-Null
